$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-24 Wednesday" "2025-09-25 Thursday"

Replace-Text "793×2=" "387×2="
Replace-Text "922×5=" "736×3="
Replace-Text "977×7=" "682×9="
Replace-Text "585×5=" "805×6="
Replace-Text "199×5=" "999×9="
Replace-Text "969×3=" "191×4="
Replace-Text "359×9=" "760×6="
Replace-Text "501×8=" "117×8="
Replace-Text "635×6=" "149×4="
Replace-Text "458×7=" "870×9="
Replace-Text "580×8=" "889×4="
Replace-Text "336×9=" "358×5="
Replace-Text "893×4=" "487×6="
Replace-Text "557×4=" "247×4="
Replace-Text "551×7=" "431×7="
Replace-Text "936×6=" "568×9="
Replace-Text "406×5=" "561×9="
Replace-Text "207×4=" "141×8="
Replace-Text "979×6=" "706×2="
Replace-Text "406×4=" "858×8="
Replace-Text "904×3=" "231×3="
Replace-Text "690×4=" "306×4="
Replace-Text "416×7=" "384×5="
Replace-Text "334×2=" "760×8="
Replace-Text "149×7=" "794×3="
